# Atualização de bases das ligas, do dia: 01-06-2024 às 01:16
#
# The source data rows for several fixtures were swapped with their
# neighbouring row (the match "id" in column A stays put, but every other
# field - match id in B, teams, scores, odds, etc. - belongs to the other
# row). This script swaps columns B:AD between each affected row pair.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose B:AD data must be swapped with each other.
$pairs = @(
    @(123, 124),
    @(128, 129),
    @(175, 176),
    @(177, 178),
    @(291, 292),
    @(296, 297),
    @(302, 304)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AD$r1")
    $range2 = $ws.Range("B$r2`:AD$r2")

    $data1 = $range1.Value2
    $data2 = $range2.Value2

    $range1.Value = $data2
    $range2.Value = $data1
}
